# Applies the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel auto-converts a plain numeric-looking string ("1.00", "0.999", ...)
# into a real number as soon as it is assigned via .Value, which would lose
# the original text formatting (trailing zeros, etc). Prefixing the string
# with a leading apostrophe -- exactly like typing it by hand into the Excel
# UI -- forces it to be stored as text instead. Values that already contain
# a second "." (e.g. "64.266.24") are never auto-parsed as numbers, so the
# prefix is only added when the value actually looks like a plain decimal.
function Set-CellText($cellRef, $text) {
    if ($text -match '^-?\d+(\.\d+)?$') {
        $ws.Range($cellRef).Value = "'" + $text
    } else {
        $ws.Range($cellRef).Value = $text
    }
}

Set-CellText 'D2' '64.266.24'
Set-CellText 'E2' '  +1.50%  '
Set-CellText 'D3' '2.648.74'
Set-CellText 'E3' '  +0.04%  '
Set-CellText 'D4' '0.999'
Set-CellText 'E4' '  -0.07%  '
Set-CellText 'D5' '604.60'
Set-CellText 'E5' '  -0.40%  '
Set-CellText 'D6' '152.91'
Set-CellText 'E6' '  +5.61%  '
Set-CellText 'E7' '  -0.03%  '
Set-CellText 'E8' '  +0.96%  '
Set-CellText 'E9' '  +1.79%  '
Set-CellText 'D10' '0.390'
Set-CellText 'E10' '  +7.40%  '
Set-CellText 'D11' '5.60'
Set-CellText 'E11' '  -0.71%  '
Set-CellText 'E12' '  -0.85%  '
Set-CellText 'D13' '27.91'
Set-CellText 'E13' '  +2.13%  '
Set-CellText 'D14' '3.120.42'
Set-CellText 'E14' '  -0.15%  '
Set-CellText 'D15' '64.038.65'
Set-CellText 'E15' '  +1.38%  '
Set-CellText 'E16' '  +1.88%  '
Set-CellText 'D17' '2.666.91'
Set-CellText 'E17' '  +1.24%  '
Set-CellText 'D18' '12.19'
Set-CellText 'E18' '  +6.91%  '
Set-CellText 'D19' '4.65'
Set-CellText 'E19' '  +4.18%  '
Set-CellText 'D20' '349.11'
Set-CellText 'E20' '  +1.97%  '
Set-CellText 'D21' '6.94'
Set-CellText 'E21' '  +1.40%  '
Set-CellText 'D22' '1.00'
Set-CellText 'E22' '  +0.08%  '
Set-CellText 'E23' '  -0.03%  '
Set-CellText 'D24' '66.62'
Set-CellText 'E24' '  -0.69%  '
Set-CellText 'D25' '1.74'
Set-CellText 'E25' '  +11.46%  '
Set-CellText 'D26' '9.44'
Set-CellText 'E26' '  +8.09%  '
Set-CellText 'D27' '1.71'
Set-CellText 'E27' '  +5.03%  '
Set-CellText 'B28' 'Bittensor'
Set-CellText 'C28' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-CellText 'D28' '559.39'
Set-CellText 'E28' '  +3.80%  '
Set-CellText 'E29' '  +0.57%  '
Set-CellText 'B30' 'Aptos'
Set-CellText 'C30' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-CellText 'D30' '8.14'
Set-CellText 'E30' '  +2.55%  '
Set-CellText 'E31' '  +0.04%  '
Set-CellText 'E32' '  +0.36%  '
Set-CellText 'D33' '0.0₃0854'
Set-CellText 'E33' '  +5.19%  '
Set-CellText 'E34' '  -0.31%  '
Set-CellText 'D35' '5.36'
Set-CellText 'E35' '  +3.22%  '
Set-CellText 'D36' '169.18'
Set-CellText 'E36' '  -1.42%  '
Set-CellText 'D37' '0.410'
Set-CellText 'E37' '  +0.82%  '
Set-CellText 'E38' '  -0.17%  '
Set-CellText 'D39' '1.97'
Set-CellText 'E39' '  +6.86%  '
Set-CellText 'D40' '19.46'
Set-CellText 'E40' '  +1.46%  '
Set-CellText 'D42' '165.23'
Set-CellText 'E42' '  -4.27%  '
Set-CellText 'D43' '40.10'
Set-CellText 'E43' '  +0.08%  '
Set-CellText 'D44' '3.88'
Set-CellText 'E44' '  +2.74%  '
Set-CellText 'D45' '0.0576'
Set-CellText 'E45' '  +0.17%  '
Set-CellText 'D46' '21.69'
Set-CellText 'E46' '  -3.47%  '
Set-CellText 'D47' '0.631'
Set-CellText 'E47' '  +0.25%  '
Set-CellText 'D48' '0.0249'
Set-CellText 'E48' '  +3.49%  '
Set-CellText 'D49' '1.98'
Set-CellText 'E49' '  +12.36%  '
Set-CellText 'D50' '0.0969'
Set-CellText 'E50' '  +0.44%  '
Set-CellText 'D51' '19.21'
Set-CellText 'E51' '  +2.15%  '
